$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: copy number/text formatting from the last filled row (18) ---
$ws.Range("B18:D18").Copy() | Out-Null
$ws.Range("B19:D19").PasteSpecial(-4122) | Out-Null
$ws.Range("F18:H18").Copy() | Out-Null
$ws.Range("F19:H19").PasteSpecial(-4122) | Out-Null

# --- Row 20: same formatting copy ---
$ws.Range("B18:D18").Copy() | Out-Null
$ws.Range("B20:D20").PasteSpecial(-4122) | Out-Null
$ws.Range("F18:H18").Copy() | Out-Null
$ws.Range("F20:H20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row 19 values (left block then right block, so new shared strings are appended in order) ---
$ws.Cells.Item(19, 2).Value = 45451
$ws.Cells.Item(19, 3).Value = "30m"
$ws.Cells.Item(19, 4).Value = "Essai pour faire fonctionner le séquenceur"

$ws.Cells.Item(19, 6).Value = 45451
$ws.Cells.Item(19, 7).Value = "30m"
$ws.Cells.Item(19, 8).Value = "Essai pour faire fonctionner le séquenceur"

# --- Row 20 values ---
$ws.Cells.Item(20, 2).Value = 45453
$ws.Cells.Item(20, 3).Value = "1h"
$ws.Cells.Item(20, 4).Value = "Continuation de testbench et des corrections"

$ws.Cells.Item(20, 6).Value = 45453
$ws.Cells.Item(20, 7).Value = "1h"
$ws.Cells.Item(20, 8).Value = "Continuation de testbench et des corrections"

# --- Restore the active selection to H16 ---
$ws.Range("H16").Select()
